$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "['HMV23/000055/0923/17', 'HMV23/000059/0923/3']"
$ws.Range("C2").Value = "['DURING ARRIVAL INSPECTION FOUND STATIC DISCHARGERS DAMAGED AT LOCATION LH WING TIP TRAILING EDGE TYPE B DISCHARGER QTY#1, RH HORZ STAB TYPE A DISCHARGER QTY#1 TYPE B DISCHARGER QTY#1, LH HORZ STAB TYPE A DISCHARGER QTY#1 TYPE B DISCHARGER QTY#3, VERT STAB TYPE A DISCHARGER QTY#1, RUDDER TYPE B DISCHARGER QTY#1.', 'DURING ARRIVAL INSPECTION FOUND STATIC DISCHARGERS DAMAGED AT LOCATION : \n1) LH HORZ STAB TYPE A DISCHARGER QTY# 2 & TYPE B DISCHARGER QTY#3 \n2) RH HORZ STAB TYPE A DISCHARGER QTY#1 \n3)VERTICAL STABILISER TYPE A DISCHARGER QTY 01 & TYPE B DISCHARGER QTY 02.']"
$ws.Range("D2").Value = "[7.0, 5.0]"

$ws.Range("B3").Value = "['HMV23/000059/0923/4', 'HMV23/000055/0923/21']"
$ws.Range("C3").Value = "['DURING ARRIVAL CHECK FOUND SERVICE LIGHT LAMP INOP AT FOLLOWING LOCATION \n1.103VU BATTERY COMPARMENT LAMP FIN: 8LS QTY#1,.\n2. REFUEL/DEFUEL CONTROL PANEL FLOODLIGHT LAMP FIN: 32QU.', 'DURING ARRIVAL CHECK FOUND SERVICE LIGHT LAMP INOP AT FOLLOWING LOCATION\n1. 80VU RH AVIONICS COMPARTMENT FIN: 9LS QTY#1, 90VU FWD AVIONICS COMPARMENT FIN: 6LS QTY#1.\n2. THS COMPARMENT FIN: 5LJ QTY#1.']"
$ws.Range("D3").Value = "[3.0, 3.0]"

$ws.Range("B4").Value = "['HMV23/000059/0923/5', 'HMV23/000055/0923/2', 'HMV23/000043/0823/12', 'HMV23/000055/0923/6']"
$ws.Range("C4").Value = "['DURING ARRIVAL INSPECTION , OBSERVED FOLLOWINS DICREPANCIES.\n\n1) RH SIDE WING NO#1 SLAT, NO.3 TRACK SQURE SEAL UPPER BULB SEAL DEGRADED.\n2) RH SIIDE WING SLAT NO#5 OUT BOARD WEATHER SEAL ERODED. SAME TO BE REPLACED.', 'DURING INSPECTION OBSERVED LH WING SLAT#5 WEATHER SEAL ERODED.SAME TO BE REPLACED.', 'WHILE ARRIVAL INSPECTION OBSERVED RH WING SLAT#5 OUTBOARD WEATHER SEAL ERODED.SAME TO BE REPLACED.', 'DURING INSPECTION, OBSERVED RH SIDE WING SLAT NO#5 WEATHER SEAL ERODED.SAME TO BE REPLACED.']"
$ws.Range("D4").Value = "[0.0, 2.0, 2.0, 3.0]"

$ws.Range("B5").Value = "['HMV23/000059/0923/6', 'HMV23/000055/0923/11']"
$ws.Range("C5").Value = "['DURING INSPECTION FOUND FOLLOWING FORWARD CARGO FLOOR PANEL DAMAGED , TO BE FABRICATED/REPLACED \n\n\n\n131EF P/N D5367521000000 \n\n131FF P/N D5367400800000 \n\n131BF P/N D5367401200000 \n\n131JF P/N D5367400600000', 'DURING INSPECTION FOUND FOLLOWING FORWARD CARGO FLOOR PANEL DAMAGED , TO BE FABRICATED/REPLACED \n\n131 AF P/N D5367520400000\n\n131EF P/N D5367521000000\n\n131FF P/N D5367400800000\n\n131BF P/N D5367401200000\n\n131CF P/N D5367401100000\n\n131GF P/N D5367400700000']"
$ws.Range("D5").Value = "[15.0, 15.0]"

$ws.Range("B7").Value = "['HMV23/000059/0923/8', 'HMV23/000055/0923/14', 'HMV23/000055/0923/13', 'HMV23/000059/0923/9']"
$ws.Range("C7").Value = "['DURING ARRIVAL INSPECTION OF FORWARD CARGO , FOUND HARDWARE DAMAGED/MISSING', 'DURING ARRIVAL INSPECTION OF AFT CARGO , FOUND \n1) CARGO SECTION HOSE MISSING/ DAMAGED \n2) FOUND HARDWARE DAMAGED/MISSING', 'DURING ARRIVAL INSPECTION OF FORWARD CARGO , FOUND \n1) CARGO SECTION HOSE MISSING/ DAMAGED \n2) FOUND HARDWARE DAMAGED/MISSING', 'DURING ARRIVAL INSPECTION OF AFT CARGO , FOUND \n1) CARGO SECTION HOSE MISSING/ DAMAGED \n2) FOUND HARDWARE DAMAGED/MISSING']"
$ws.Range("D7").Value = "[5.0, 5.0, 5.0, 7.0]"

$ws.Range("B11").Value = "['HMV23/000055/0923/10', 'HMV23/000059/0923/15']"
$ws.Range("C11").Value = "['DURING WEEKLY INSPECTION OBSERVED FIN: 11WL EPSU BATTERY INOP.', 'DURING ARRIVAL CHECK FOUND EPSU BATTERIES QTY#2 FIN:10WL AND FIN: 11WL INOP.']"
$ws.Range("D11").Value = "[5.0, 6.0]"

$ws.Range("B12").Value = "['HMV23/000055/0923/23', 'HMV23/000059/0923/16']"
$ws.Range("C12").Value = "['DURING WEAR LIMIT BONDING LEADS INSPECTION CHECK FOUND BONDING LEADS DAMAGED AT LOCATION TOILET SERVICE PANEL QTY#2, WATER SERVICE PANEL QTY#2, RH MLG DOOR QTY#1, LH MLG DOOR QTY#2, 108VU QTY#1, FWD CARGO DOOR QTY#1, AFT CARGO DOOR QTY#1, FWD LH PASSENGER DOOR QTY#1, RH WING FLAP TRACK FAIRING NO.2 QTY#1, RH WING FLAP TRACK FAIRING NO.3 QTY#1, LH WING FLAP TRACK FAIRING NO.1 QTY#2, LH WING FLAP TRACK FAIRING NO.3 QTY#2, LH ELEVATOR PANEL 334BB QTY#1, LH ELEVATOR PANEL 334HB QTY#1, RH ELEVATOR PANEL 344DB QTY#1, RH ELEVATOR PANEL 344BB QTY#1.', 'DURING WEAR LIMIT BONDING LEADS INSPECTION CHECK FOUND BONDING LEADS DAMAGED AT LOCATION TOILET SERVICE PANEL QTY#1, WATER SERVICE PANEL QTY#2, RH MLG DOOR QTY#1, LH MLG DOOR QTY#1, LH ELEVATOR QTY#3, RH ELEVATOR QTY#1, RH WING BETWEEN SPOILER 3 AND 4 QTY#2, LH WING FLAP TRACK FAIRING NO.3 QTY#1, RH WING FLAP TRACK FAIRING NO.1 QTY#2.']"
$ws.Range("D12").Value = "[24.0, 18.0]"

$ws.Range("B14").Value = "['HMV23/000055/0923/1', 'HMV23/000055/0923/7', 'HMV23/000055/0923/9', 'HMV23/000055/0923/8']"
$ws.Range("C14").Value = "['DURING ARRIVAL INSPECTION , OBSERVED DENT ON SLAT NO #2 LIP AREA AT RH SIDE WING. ASSESSMENT TO BE CARRY OUT.\n\nSLAT NO #2 P/N- D5746092000500\n SA7205698', 'DURING ARRIVAL INSPECTION , OBSERVED DENT ON SLAT NO #2 LIP AREA AT RH SIDE WING. DENT LOCATION, 3 INCH FROM T/E AND 67.5 INCH FROM OUTBD EDGE .ASSESSMENT TO BE CARRY OUT.\nSLAT NO #2 P/N- D5746092000500\n S/N- SA7205698', 'DURING ARRIVAL INSPECTION, OBSERVED SMALL DENTS ON SLAT NO #1 TOP SKIN AT RH SIDE WING. LOCATION OF DENT 3.5 INCH FROM T/E AND 37 INCH FROM INBD EDGE OF SLAT.\nSLAT NO#1 P/N- D5746091000300\n S/N- SA7205917', 'DURING ARRIVAL INSPECTION, OBSERVED SMALL DENTS ON SLAT NO #1 TOP SKIN AT RH SIDE WING. LOCATION OF DENT 3.5 INCH FROM T/E AND 37 INCH FROM INBD EDGE OF SLAT.\nSLAT NO#1 P/N- D5746091000300\n S/N- SA7205917']"
$ws.Range("D14").Value = "[2.0, 7.0, 7.0, 3.0]"

$ws.Range("B15").Value = "['HMV23/000055/0923/4', 'HMV23/000055/0923/3', 'HMV23/000055/0923/5']"
$ws.Range("C15").Value = "['DURING ARRIVAL INSPECTION FOUND FOLLOWING FAILURE MSG IN PFR:\n1. 34-53-31 ADF1(2RP1)\nFURTHER RECTIFICATION TO BE CARRIED OUT.', 'DURING ARRIVAL INSPECTION FOUND FOLLOWING WARNING MSG IN PFR:\n1. AUTO FLT AP OFF\nFURTHER RECTIFICATION TO BE CARRIED OUT.', 'DURING ARRIVAL INSPECTION FOUND FOLLOWING FAILURE MSG IN PFR:\n1. 38-31-41 TOILET ASSY LAV F\nFURTHER RECTIFICATION TO BE CARRIED OUT.']"
$ws.Range("D15").Value = "[2.5, 2.0, 4.0]"
